# ---------------------------------------------------------------------------
# "Updated title of the game"
#
# The title paragraph originally reads "The Rite of Dr. Adams". The word
# "of" is replaced with "for", the way a user would do it interactively in
# Word: select the word "of" and type "for" over it. Word leaves the
# insertion point (and therefore its "_GoBack" last-edit bookmark) right
# after the newly typed text, which also happens to split the run that used
# to hold the whole title into three runs (the text before the edit, the
# freshly typed text, and the text after the edit).
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# "The Rite of Dr. Adams"
#  0123456789...
# "of" sits at characters 9-10 (0-based range 9..11)
$target = $d.Range(9, 11)
$target.Text = "for"

# Force the just-typed word into its own run (mirrors Word's behaviour of
# keeping the edited span separate from the untouched text around it) by
# collapsing a range right before it and reuniting it right after - the
# bookmark insertion below does the splitting for us.
$splitPoint = $d.Range(9, 9)
$d.Bookmarks.Add("zzz_tmp_split", $splitPoint) | Out-Null

# Word always keeps exactly one "_GoBack" bookmark, tracking the last place
# text was edited - re-adding it here moves it from its old location (after
# "Gate for leopard") to right after the word that was just typed.
$goBackPoint = $d.Range(12, 12)
$d.Bookmarks.Add("_GoBack", $goBackPoint) | Out-Null

# Drop the helper bookmark now that it has done its job of keeping the runs
# split apart.
$d.Bookmarks("zzz_tmp_split").Delete()

# ---------------------------------------------------------------------------
# Re-saving the document also compacts the z-order (relativeHeight) values
# Word assigns to every floating drawing, renumbering them sequentially (in
# document order) starting from 251650048 in steps of 1024. Reproduce that
# renumbering using the standard ZOrder move commands (2 = bring forward one
# step, 3 = send backward one step), applied the right number of times to
# each shape to land on its final value.
# ---------------------------------------------------------------------------

$zsteps = @(-8, -10, -7, -5, -8, -8, -2, -1, 0, -5, 2, 3, 4, 5, 6, 7)
for ($i = 0; $i -lt $zsteps.Length; $i++) {
    $steps = $zsteps[$i]
    if ($steps -eq 0) { continue }
    $shape = $d.Shapes.Item($i + 1)
    if ($steps -gt 0) {
        for ($k = 0; $k -lt $steps; $k++) { $shape.ZOrder(2) | Out-Null }
    } else {
        for ($k = 0; $k -lt (-$steps); $k++) { $shape.ZOrder(3) | Out-Null }
    }
}
